$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 613.0769
$ws.Range("I2").Value = 130.11111
$ws.Range("J2").Value = 1699.75
$ws.Range("K2").Value = 130.11111
$ws.Range("L2").Value = 1699.75
$ws.Range("M2").Value = -17.11111
$ws.Range("N2").Value = -1925.75
$ws.Range("H40").Value = 3500.5
$ws.Range("I40").Value = 5999.6665
$ws.Range("J40").Value = 1001.3333
$ws.Range("K40").Value = 5999.6665
$ws.Range("L40").Value = 1001.3333
$ws.Range("M40").Value = -5824.6665
$ws.Range("N40").Value = -1351.3333
$ws.Range("H42").Value = 2104.3
$ws.Range("I42").Value = 637.375
$ws.Range("K42").Value = 1912.125
$ws.Range("M42").Value = -1682.125
$ws.Range("H69").Value = 17993.438
$ws.Range("J69").Value = 20000
$ws.Range("L69").Value = 60000
$ws.Range("N69").Value = -61748
$ws.Range("H72").Value = 17993.438
$ws.Range("J72").Value = 20000
$ws.Range("L72").Value = 180000
$ws.Range("N72").Value = -188736
$ws.Range("H98").Value = 2708.7297
$ws.Range("I98").Value = 808.5357
$ws.Range("J98").Value = 8620.444
$ws.Range("K98").Value = 808.5357
$ws.Range("L98").Value = 8620.444
$ws.Range("M98").Value = 689.4643
$ws.Range("N98").Value = -11616.444
$ws.Range("H112").Value = 3892.2903
$ws.Range("J112").Value = 3976.7334
$ws.Range("L112").Value = 11930.2002
$ws.Range("N112").Value = -14146.2002
$ws.Range("H122").Value = 2708.7297
$ws.Range("I122").Value = 808.5357
$ws.Range("J122").Value = 8620.444
$ws.Range("K122").Value = 2425.6071
$ws.Range("L122").Value = 25861.332
$ws.Range("M122").Value = 24.39289999999983
$ws.Range("N122").Value = -30761.332
$ws.Range("H130").Value = 94963.336
$ws.Range("J130").Value = 94963.336
$ws.Range("L130").Value = 94963.336
$ws.Range("N130").Value = -105003.336
$ws.Range("H138").Value = 2853.0571
$ws.Range("J138").Value = 5262.533
$ws.Range("L138").Value = 15787.599
$ws.Range("N138").Value = -26067.599

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4263.518
$ws.Range("I32").Value = 563.65955
$ws.Range("K32").Value = 563.65955
$ws.Range("M32").Value = -276.65955
$ws.Range("H43").Value = 15821.375
$ws.Range("I43").Value = 27447.5
$ws.Range("J43").Value = 11946
$ws.Range("K43").Value = 27447.5
$ws.Range("L43").Value = 11946
$ws.Range("M43").Value = -27134.5
$ws.Range("N43").Value = -12572
$ws.Range("H61").Value = 5690.25
$ws.Range("J61").Value = 13194
$ws.Range("L61").Value = 13194
$ws.Range("N61").Value = -13618
$ws.Range("H132").Value = 3195.394
$ws.Range("I132").Value = 2833.6428
$ws.Range("J132").Value = 5221.2
$ws.Range("K132").Value = 8500.928400000001
$ws.Range("L132").Value = 15663.6
$ws.Range("M132").Value = -5970.928400000001
$ws.Range("N132").Value = -20723.6
$ws.Range("H136").Value = 5690.25
$ws.Range("J136").Value = 13194
$ws.Range("L136").Value = 39582
$ws.Range("N136").Value = -44682

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 299999
$ws.Range("J55").Value = 299999
$ws.Range("L55").Value = 299999
$ws.Range("N55").Value = -300545

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 20992.732
$ws.Range("I16").Value = 15317.546
$ws.Range("K16").Value = 15317.546
$ws.Range("M16").Value = -15030.546
$ws.Range("H99").Value = 6791
$ws.Range("I99").Value = 7400.5713
$ws.Range("K99").Value = 7400.5713
$ws.Range("M99").Value = -5902.5713
$ws.Range("H113").Value = 20992.732
$ws.Range("I113").Value = 15317.546
$ws.Range("K113").Value = 15317.546
$ws.Range("M113").Value = -13147.546
$ws.Range("H126").Value = 6791
$ws.Range("I126").Value = 7400.5713
$ws.Range("K126").Value = 22201.7139
$ws.Range("M126").Value = -19731.7139

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 557.6
$ws.Range("I5").Value = 411.22726
$ws.Range("K5").Value = 1233.68178
$ws.Range("M5").Value = -1121.68178
$ws.Range("H38").Value = 81.181816
$ws.Range("J38").Value = 105.57143
$ws.Range("L38").Value = 316.71429
$ws.Range("N38").Value = -1010.71429
$ws.Range("H41").Value = 735.2941
$ws.Range("J41").Value = 500
$ws.Range("L41").Value = 1500
$ws.Range("N41").Value = -2176
$ws.Range("H55").Value = 5164.2085
$ws.Range("I55").Value = 986.1818
$ws.Range("J55").Value = 8699.462
$ws.Range("K55").Value = 2958.5454
$ws.Range("L55").Value = 26098.386
$ws.Range("M55").Value = -2781.5454
$ws.Range("N55").Value = -26452.386
$ws.Range("H58").Value = 8600.666999999999
$ws.Range("J58").Value = 8600.666999999999
$ws.Range("L58").Value = 25802.001
$ws.Range("N58").Value = -26058.001
$ws.Range("H80").Value = 2999
$ws.Range("J80").Value = 2999
$ws.Range("L80").Value = 8997
$ws.Range("N80").Value = -10869
$ws.Range("H83").Value = 2999
$ws.Range("J83").Value = 2999
$ws.Range("L83").Value = 26991
$ws.Range("N83").Value = -36351
$ws.Range("H113").Value = 2406.077
$ws.Range("J113").Value = 2527.4167
$ws.Range("L113").Value = 7582.250100000001
$ws.Range("N113").Value = -11922.2501
$ws.Range("H119").Value = 6053.5713
$ws.Range("I119").Value = 3729.3333
$ws.Range("K119").Value = 11187.9999
$ws.Range("M119").Value = -6349.999899999999
$ws.Range("H123").Value = 9384.117
$ws.Range("I123").Value = 4765
$ws.Range("J123").Value = 10000
$ws.Range("K123").Value = 14295
$ws.Range("L123").Value = 30000
$ws.Range("M123").Value = -11845
$ws.Range("N123").Value = -34900
$ws.Range("H131").Value = 1920.4286
$ws.Range("J131").Value = 2613.2222
$ws.Range("L131").Value = 7839.6666
$ws.Range("N131").Value = -17919.6666
$ws.Range("H135").Value = 557.6
$ws.Range("I135").Value = 411.22726
$ws.Range("K135").Value = 3701.04534
$ws.Range("M135").Value = -1166.04534
$ws.Range("H140").Value = 1085.909
$ws.Range("I140").Value = 1026.0938
$ws.Range("K140").Value = 3078.2814
$ws.Range("M140").Value = 2101.7186

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2148.3333
$ws.Range("I102").Value = 1709.6923
$ws.Range("K102").Value = 1709.6923
$ws.Range("M102").Value = -87.69229999999993
$ws.Range("H122").Value = 3375.7778
$ws.Range("I122").Value = 3297.8125
$ws.Range("K122").Value = 9893.4375
$ws.Range("M122").Value = -7443.4375
$ws.Range("H132").Value = 5495.1333
$ws.Range("I132").Value = 5858.636
$ws.Range("K132").Value = 17575.908
$ws.Range("M132").Value = -15045.908

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2123.95
$ws.Range("I136").Value = 1792.8823
$ws.Range("K136").Value = 5378.6469
$ws.Range("M136").Value = -2828.6469

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 99989.664
$ws.Range("J45").Value = 99989.664
$ws.Range("L45").Value = 99989.664
$ws.Range("N45").Value = -100971.664
$ws.Range("H57").Value = 105000
$ws.Range("J57").Value = 105000
$ws.Range("L57").Value = 105000
$ws.Range("N57").Value = -106508
$ws.Range("I132").Value = 2582.6135
$ws.Range("J132").Value = 2918
$ws.Range("K132").Value = 7747.8405
$ws.Range("L132").Value = 8754
$ws.Range("M132").Value = -5217.8405
$ws.Range("N132").Value = -13814
$ws.Range("H136").Value = 3357.8276
$ws.Range("I136").Value = 3334.8928
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 10004.6784
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -7454.678400000001
$ws.Range("N136").Value = -17100
